# Fix error handling for bad file extensions
# Update time-series snippet data rows (cols A-E) to reflect corrected snippet extraction
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 182
$ws.Range("B7").Value = "6:56 AM"
$ws.Range("C7").Value = 5440.89
$ws.Range("D7").Value = 5470.89
$ws.Range("E7").Value = 108.74

$ws.Range("A8").Value = 230
$ws.Range("B8").Value = "7:20 AM"
$ws.Range("C8").Value = 6880.786667
$ws.Range("D8").Value = 6910.786667
$ws.Range("E8").Value = 32.193333

$ws.Range("A9").Value = 239
$ws.Range("B9").Value = "7:24 AM"
$ws.Range("C9").Value = 7146.62
$ws.Range("D9").Value = 7176.62
$ws.Range("E9").Value = 32.95

$ws.Range("A10").Value = 251
$ws.Range("B10").Value = "7:30 AM"
$ws.Range("C10").Value = 7516.66
$ws.Range("D10").Value = 7546.66
$ws.Range("E10").Value = 27.63

$ws.Range("A13").Value = 550
$ws.Range("B13").Value = "10:00 A"
$ws.Range("C13").Value = 16484
$ws.Range("D13").Value = 16514
$ws.Range("E13").Value = 26.08

$ws.Range("A14").Value = 570
$ws.Range("B14").Value = "10:10 A"
$ws.Range("C14").Value = 17087.03
$ws.Range("D14").Value = 17117.03
$ws.Range("E14").Value = 37.05

$ws.Range("A15").Value = 993
$ws.Range("B15").Value = "1:41 PM"
$ws.Range("C15").Value = 29775.52
$ws.Range("D15").Value = 29805.52
$ws.Range("E15").Value = 41.35

$ws.Range("A16").Value = 1016
$ws.Range("B16").Value = "1:53 PM"
$ws.Range("C16").Value = 30454.73
$ws.Range("D16").Value = 30484.73
$ws.Range("E16").Value = 24.88

$ws.Range("A17").Value = 1033
$ws.Range("B17").Value = "2:01 PM"
$ws.Range("C17").Value = 30974.42
$ws.Range("D17").Value = 31004.42
$ws.Range("E17").Value = 36.19

$ws.Range("A21").Value = 1180
$ws.Range("B21").Value = "3:15 PM"
$ws.Range("C21").Value = 35378.663333
$ws.Range("D21").Value = 35408.663333
$ws.Range("E21").Value = 42.293333

$ws.Range("A22").Value = 1190
$ws.Range("B22").Value = "3:20 PM"
$ws.Range("C22").Value = 35678.675
$ws.Range("D22").Value = 35708.675
$ws.Range("E22").Value = 40.53

$ws.Range("A23").Value = 1196
$ws.Range("B23").Value = "3:23 PM"
$ws.Range("C23").Value = 35858.715
$ws.Range("D23").Value = 35888.715
$ws.Range("E23").Value = 56.205

$ws.Range("A24").Value = 1202
$ws.Range("B24").Value = "3:26 PM"
$ws.Range("C24").Value = 36041.01
$ws.Range("D24").Value = 36071.01
$ws.Range("E24").Value = 24.126667

$ws.Range("A25").Value = 1209
$ws.Range("B25").Value = "3:29 PM"
$ws.Range("C25").Value = 36244.02
$ws.Range("D25").Value = 36274.02
$ws.Range("E25").Value = 75.29000000000001

$ws.Range("A26").Value = 1221
$ws.Range("B26").Value = "3:35 PM"
$ws.Range("C26").Value = 36616.4
$ws.Range("D26").Value = 36646.4
$ws.Range("E26").Value = 23.47

$ws.Range("A27").Value = 1227
$ws.Range("B27").Value = "3:39 PM"
$ws.Range("C27").Value = 36808.87
$ws.Range("D27").Value = 36838.87
$ws.Range("E27").Value = 22.12

$ws.Range("A28").Value = 1245
$ws.Range("B28").Value = "3:47 PM"
$ws.Range("C28").Value = 37332.09
$ws.Range("D28").Value = 37362.09
$ws.Range("E28").Value = 193.12

$ws.Range("A29").Value = 1257
$ws.Range("B29").Value = "3:53 PM"
$ws.Range("C29").Value = 37684.8
$ws.Range("D29").Value = 37714.8
$ws.Range("E29").Value = 36.82

$ws.Range("A30").Value = 1269
$ws.Range("B30").Value = "3:59 PM"
$ws.Range("C30").Value = 38052.733333
$ws.Range("D30").Value = 38082.733333
$ws.Range("E30").Value = 20.703333

$ws.Range("A31").Value = 1278
$ws.Range("B31").Value = "4:04 PM"
$ws.Range("C31").Value = 38317.495
$ws.Range("D31").Value = 38347.495
$ws.Range("E31").Value = 26.47

$ws.Range("A32").Value = 1288
$ws.Range("B32").Value = "4:09 PM"
$ws.Range("C32").Value = 38622.23
$ws.Range("D32").Value = 38652.23
$ws.Range("E32").Value = 50.7

$ws.Range("A33").Value = 1302
$ws.Range("B33").Value = "4:16 PM"
$ws.Range("C33").Value = 39032.49
$ws.Range("D33").Value = 39062.49
$ws.Range("E33").Value = 230.5

$ws.Range("A34").Value = 1310
$ws.Range("B34").Value = "4:20 PM"
$ws.Range("C34").Value = 39277.12
$ws.Range("D34").Value = 39307.12
$ws.Range("E34").Value = 28.49

$ws.Range("A35").Value = 1323
$ws.Range("B35").Value = "4:26 PM"
$ws.Range("C35").Value = 39677.255
$ws.Range("D35").Value = 39707.255
$ws.Range("E35").Value = 57.735

$ws.Range("A36").Value = 1333
$ws.Range("B36").Value = "4:31 PM"
$ws.Range("C36").Value = 39975.185
$ws.Range("D36").Value = 40005.185
$ws.Range("E36").Value = 26.385

$ws.Range("A37").Value = 1341
$ws.Range("B37").Value = "4:35 PM"
$ws.Range("C37").Value = 40209.895
$ws.Range("D37").Value = 40239.895
$ws.Range("E37").Value = 54.99

$ws.Range("A38").Value = 1349
$ws.Range("B38").Value = "4:39 PM"
$ws.Range("C38").Value = 40446.58
$ws.Range("D38").Value = 40476.58
$ws.Range("E38").Value = 52.82

$ws.Range("A39").Value = 1360
$ws.Range("B39").Value = "4:45 PM"
$ws.Range("C39").Value = 40775.5
$ws.Range("D39").Value = 40805.5
$ws.Range("E39").Value = 27.96

$ws.Range("A40").Value = 1370
$ws.Range("B40").Value = "4:50 PM"
$ws.Range("C40").Value = 41074.77
$ws.Range("D40").Value = 41104.77
$ws.Range("E40").Value = 24.045

$ws.Range("A41").Value = 1390
$ws.Range("B41").Value = "5:00 PM"
$ws.Range("C41").Value = 41680.48
$ws.Range("D41").Value = 41710.48
$ws.Range("E41").Value = 38.76

$ws.Range("A42").Value = 1410
$ws.Range("B42").Value = "5:10 PM"
$ws.Range("C42").Value = 42278.335
$ws.Range("D42").Value = 42308.335
$ws.Range("E42").Value = 44.55

$ws.Range("A43").Value = 1417
$ws.Range("B43").Value = "5:13 PM"
$ws.Range("C43").Value = 42486.06
$ws.Range("D43").Value = 42516.06
$ws.Range("E43").Value = 61.41

$ws.Range("A44").Value = 1427
$ws.Range("B44").Value = "5:19 PM"
$ws.Range("C44").Value = 42809.97
$ws.Range("D44").Value = 42839.97
$ws.Range("E44").Value = 161.82

$ws.Range("A45").Value = 1441
$ws.Range("B45").Value = "5:26 PM"
$ws.Range("C45").Value = 43229.93
$ws.Range("D45").Value = 43259.93
$ws.Range("E45").Value = 24.11

$ws.Range("A46").Value = 1451
$ws.Range("B46").Value = "5:30 PM"
$ws.Range("C46").Value = 43516.02
$ws.Range("D46").Value = 43546.02
$ws.Range("E46").Value = 88.58

$ws.Range("A47").Value = 1472
$ws.Range("B47").Value = "5:41 PM"
$ws.Range("C47").Value = 44138.865
$ws.Range("D47").Value = 44168.865
$ws.Range("E47").Value = 48.69

$ws.Range("A48").Value = 1479
$ws.Range("B48").Value = "5:44 PM"
$ws.Range("C48").Value = 44340.68
$ws.Range("D48").Value = 44370.68
$ws.Range("E48").Value = 49.65

$ws.Range("A49").Value = 1492
$ws.Range("B49").Value = "5:51 PM"
$ws.Range("C49").Value = 44738.87
$ws.Range("D49").Value = 44768.87
$ws.Range("E49").Value = 22.62

$ws.Range("A50").Value = 1516
$ws.Range("B50").Value = "6:03 PM"
$ws.Range("C50").Value = 45455.86
$ws.Range("D50").Value = 45485.86
$ws.Range("E50").Value = 28.36

$ws.Range("A51").Value = 1551
$ws.Range("B51").Value = "6:20 PM"
$ws.Range("C51").Value = 46502.28
$ws.Range("D51").Value = 46532.28
$ws.Range("E51").Value = 39.27
